$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect the "Price" column cells we are about to rewrite so Excel keeps
# them as text (matching the original inline-string cells) instead of
# auto-converting number-like strings (e.g. "88.50") into numeric values.
$priceRows = @(2,3,5,6,10,11,12,15,17,19,20,23,24,28,31,32,34,35,37,40,42,46,47,49,51)
foreach ($r in $priceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

# Price (column D) updates
$ws.Range("D2").Value = "40.021.12"
$ws.Range("D3").Value = "2.217.21"
$ws.Range("D5").Value = "290.37"
$ws.Range("D6").Value = "88.50"
$ws.Range("D10").Value = "30.79"
$ws.Range("D11").Value = "0.0783"
$ws.Range("D12").Value = "47.90"
$ws.Range("D15").Value = "2.559.59"
$ws.Range("D17").Value = "2.231.04"
$ws.Range("D19").Value = "39.962.91"
$ws.Range("D20").Value = "11.88"
$ws.Range("D23").Value = "65.64"
$ws.Range("D24").Value = "235.50"
$ws.Range("D28").Value = "22.64"
$ws.Range("D31").Value = "153.13"
$ws.Range("D32").Value = "32.24"
$ws.Range("D34").Value = "4.97"
$ws.Range("D35").Value = "0.0722"
$ws.Range("D37").Value = "2.85"
$ws.Range("D40").Value = "15.86"
$ws.Range("D42").Value = "2.100.30"
$ws.Range("D46").Value = "17.72"
$ws.Range("D47").Value = "9.88"
$ws.Range("D49").Value = "2.432.30"
$ws.Range("D51").Value = "88.77"

# Volume(1h) (column E) updates
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -2.73%  "
$ws.Range("E6").Value = "  +5.95%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("E10").Value = "  +3.95%  "
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("E12").Value = "  +3.98%  "
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("E14").Value = "  +3.18%  "
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("E20").Value = "  +14.39%  "
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  +1.75%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("E31").Value = "  +2.70%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  +2.91%  "
$ws.Range("E35").Value = "  +2.93%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("E37").Value = "  +7.23%  "
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("E39").Value = "  +3.29%  "
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("E41").Value = "  +3.55%  "
$ws.Range("E42").Value = "  +8.88%  "
$ws.Range("E43").Value = "  +5.07%  "
$ws.Range("E44").Value = "  +1.86%  "
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("E46").Value = "  +9.38%  "
$ws.Range("E47").Value = "  +7.18%  "
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("E51").Value = "  +0.15%  "
